$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new week's record (2023-01-20) at the top of the
# Mango price history (rows 925-927), shifting all existing rows down by
# three (one record = Especial/Primera/Segunda). This naturally pushes the
# former last record (rows 1039-1041) down to become new rows 1042-1044.

$ws.Range("A925:A927").EntireRow.Insert()

# Row 925 - Especial
$ws.Cells.Item(925,1).Value2 = 8
$ws.Cells.Item(925,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(925,3).Value2 = "Coquimbo"
$ws.Cells.Item(925,4).Value2 = 44946
$ws.Cells.Item(925,5).Value2 = 4
$ws.Cells.Item(925,6).Value2 = "Fruta"
$ws.Cells.Item(925,7).Value2 = 100108
$ws.Cells.Item(925,8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(925,9).Value2 = 100108002
$ws.Cells.Item(925,10).Value2 = "Mango"
$ws.Cells.Item(925,11).Value2 = "Sin especificar"
$ws.Cells.Item(925,12).Value2 = "Especial"
$ws.Cells.Item(925,13).Value2 = 512
$ws.Cells.Item(925,14).Value2 = 6500
$ws.Cells.Item(925,15).Value2 = 7000
$ws.Cells.Item(925,16).Value2 = 6750
$ws.Cells.Item(925,17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(925,18).Value2 = "Perú"
$ws.Cells.Item(925,19).Value2 = 1688
$ws.Cells.Item(925,20).Value2 = 4

# Row 926 - Primera
$ws.Cells.Item(926,1).Value2 = 8
$ws.Cells.Item(926,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(926,3).Value2 = "Coquimbo"
$ws.Cells.Item(926,4).Value2 = 44946
$ws.Cells.Item(926,5).Value2 = 4
$ws.Cells.Item(926,6).Value2 = "Fruta"
$ws.Cells.Item(926,7).Value2 = 100108
$ws.Cells.Item(926,8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(926,9).Value2 = 100108002
$ws.Cells.Item(926,10).Value2 = "Mango"
$ws.Cells.Item(926,11).Value2 = "Sin especificar"
$ws.Cells.Item(926,12).Value2 = "Primera"
$ws.Cells.Item(926,13).Value2 = 512
$ws.Cells.Item(926,14).Value2 = 6500
$ws.Cells.Item(926,15).Value2 = 7000
$ws.Cells.Item(926,16).Value2 = 6750
$ws.Cells.Item(926,17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(926,18).Value2 = "Perú"
$ws.Cells.Item(926,19).Value2 = 1688
$ws.Cells.Item(926,20).Value2 = 4

# Row 927 - Segunda
$ws.Cells.Item(927,1).Value2 = 8
$ws.Cells.Item(927,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(927,3).Value2 = "Coquimbo"
$ws.Cells.Item(927,4).Value2 = 44946
$ws.Cells.Item(927,5).Value2 = 4
$ws.Cells.Item(927,6).Value2 = "Fruta"
$ws.Cells.Item(927,7).Value2 = 100108
$ws.Cells.Item(927,8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(927,9).Value2 = 100108002
$ws.Cells.Item(927,10).Value2 = "Mango"
$ws.Cells.Item(927,11).Value2 = "Sin especificar"
$ws.Cells.Item(927,12).Value2 = "Segunda"
$ws.Cells.Item(927,13).Value2 = 512
$ws.Cells.Item(927,14).Value2 = 6500
$ws.Cells.Item(927,15).Value2 = 7000
$ws.Cells.Item(927,16).Value2 = 6750
$ws.Cells.Item(927,17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(927,18).Value2 = "Perú"
$ws.Cells.Item(927,19).Value2 = 1688
$ws.Cells.Item(927,20).Value2 = 4

Write-Host "Dimension rows now:" $ws.UsedRange.Rows.Count
